$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3
$ws.Range("A3").Value = 86.5999999999986
$ws.Range("A5").Value = 24.19999999999885
$ws.Range("A6").Value = 25.19999999999885
$ws.Range("A8").Value = 30.19999999999891
$ws.Range("A9").Value = 100.5999999999985
